# Fruta / hortaliza, semanal
# Update weekly price records: the date (D), volume (M), min/max/avg price
# (N/O/P) and price per kg (S) columns are refreshed for each data row
# (rows 2-21) on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44889
$ws.Range("M2").Value = 460
$ws.Range("N2").Value = 3500
$ws.Range("O2").Value = 4000
$ws.Range("P2").Value = 3750
$ws.Range("S2").Value = 1875

$ws.Range("D3").Value = 44881
$ws.Range("M3").Value = 440
$ws.Range("N3").Value = 6000
$ws.Range("O3").Value = 7000
$ws.Range("P3").Value = 6500
$ws.Range("S3").Value = 3250

$ws.Range("D4").Value = 45211
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = 10000
$ws.Range("O4").Value = 11000
$ws.Range("P4").Value = 10500
$ws.Range("S4").Value = 5250

$ws.Range("D5").Value = 44819
$ws.Range("M5").Value = 240
$ws.Range("N5").Value = 11000
$ws.Range("O5").Value = 12000
$ws.Range("P5").Value = 11500
$ws.Range("S5").Value = 5750

$ws.Range("D6").Value = 44875
$ws.Range("M6").Value = 400
$ws.Range("N6").Value = 7000
$ws.Range("O6").Value = 7500
$ws.Range("P6").Value = 7250
$ws.Range("S6").Value = 3625

$ws.Range("D7").Value = 44455
$ws.Range("M7").Value = 200
$ws.Range("N7").Value = 12000
$ws.Range("O7").Value = 13000
$ws.Range("P7").Value = 12500
$ws.Range("S7").Value = 6250

$ws.Range("D8").Value = 44489
$ws.Range("M8").Value = 160
$ws.Range("N8").Value = 9500
$ws.Range("O8").Value = 10000
$ws.Range("P8").Value = 9750
$ws.Range("S8").Value = 4875

$ws.Range("D9").Value = 44461
$ws.Range("M9").Value = 200
$ws.Range("N9").Value = 11000
$ws.Range("O9").Value = 12000
$ws.Range("P9").Value = 11500
$ws.Range("S9").Value = 5750

$ws.Range("D10").Value = 44497
$ws.Range("M10").Value = 500
$ws.Range("N10").Value = 9000
$ws.Range("O10").Value = 10000
$ws.Range("P10").Value = 9500
$ws.Range("S10").Value = 4750

$ws.Range("D11").Value = 45204
$ws.Range("M11").Value = 400
$ws.Range("N11").Value = 9000
$ws.Range("O11").Value = 10000
$ws.Range("P11").Value = 9500
$ws.Range("S11").Value = 4750

$ws.Range("D12").Value = 44490
$ws.Range("M12").Value = 400
$ws.Range("N12").Value = 9500
$ws.Range("O12").Value = 10000
$ws.Range("P12").Value = 9750
$ws.Range("S12").Value = 4875

$ws.Range("D13").Value = 45203
$ws.Range("M13").Value = 1000
$ws.Range("N13").Value = 9000
$ws.Range("O13").Value = 10000
$ws.Range("P13").Value = 9500
$ws.Range("S13").Value = 4750

$ws.Range("D14").Value = 44874
$ws.Range("M14").Value = 300
$ws.Range("N14").Value = 7500
$ws.Range("O14").Value = 8000
$ws.Range("P14").Value = 7750
$ws.Range("S14").Value = 3875

$ws.Range("D15").Value = 44882
$ws.Range("M15").Value = 440
$ws.Range("N15").Value = 6000
$ws.Range("O15").Value = 7000
$ws.Range("P15").Value = 6500
$ws.Range("S15").Value = 3250

$ws.Range("D16").Value = 44482
$ws.Range("M16").Value = 240
$ws.Range("N16").Value = 10000
$ws.Range("O16").Value = 11000
$ws.Range("P16").Value = 10500
$ws.Range("S16").Value = 5250

$ws.Range("D17").Value = 44475
$ws.Range("M17").Value = 240
$ws.Range("N17").Value = 11000
$ws.Range("O17").Value = 12000
$ws.Range("P17").Value = 11500
$ws.Range("S17").Value = 5750

$ws.Range("D18").Value = 44454
$ws.Range("M18").Value = 160
$ws.Range("N18").Value = 12000
$ws.Range("O18").Value = 13000
$ws.Range("P18").Value = 12500
$ws.Range("S18").Value = 6250

$ws.Range("D19").Value = 44895
$ws.Range("M19").Value = 240
$ws.Range("N19").Value = 3000
$ws.Range("O19").Value = 3500
$ws.Range("P19").Value = 3250
$ws.Range("S19").Value = 1625

$ws.Range("D20").Value = 44517
$ws.Range("M20").Value = 400
$ws.Range("N20").Value = 5500
$ws.Range("O20").Value = 6000
$ws.Range("P20").Value = 5750
$ws.Range("S20").Value = 2875

$ws.Range("D21").Value = 44818
$ws.Range("M21").Value = 200
$ws.Range("N21").Value = 11000
$ws.Range("O21").Value = 12000
$ws.Range("P21").Value = 11500
$ws.Range("S21").Value = 5750
